$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the row 5 "custom accuracy" readings to 2 decimal places.
$ws.Range("B5").Value = 2.32
$ws.Range("C5").Value = 1.55
$ws.Range("D5").Value = 0.63
$ws.Range("E5").Value = 5.24
$ws.Range("F5").Value = 3.69
$ws.Range("G5").Value = 1.76
$ws.Range("H5").Value = 12.22
$ws.Range("I5").Value = 2.91
$ws.Range("J5").Value = 1.23
$ws.Range("K5").Value = 1.51
$ws.Range("L5").Value = 2.08
$ws.Range("M5").Value = 2.28
$ws.Range("N5").Value = 0.62
$ws.Range("O5").Value = 1.88
$ws.Range("P5").Value = 2.67
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 0.66
$ws.Range("S5").Value = 0.3
$ws.Range("T5").Value = 21.37
$ws.Range("U5").Value = 5.63
$ws.Range("V5").Value = 1.74
$ws.Range("W5").Value = 3.62
$ws.Range("X5").Value = 1.78
$ws.Range("Y5").Value = 0.5600000000000001
$ws.Range("Z5").Value = 5.64
$ws.Range("AA5").Value = 1.53
$ws.Range("AB5").Value = 1.5
$ws.Range("AC5").Value = 1.74
$ws.Range("AD5").Value = 2.11
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 11.74
$ws.Range("AG5").Value = 0.86
$ws.Range("AH5").Value = 2.18

# Remove the now-redundant row of raw (1000-reading) data.
$ws.Rows(6).Delete()
